$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet: insert a new (blank) column before column N,
# pushing the existing "Late" / "heading" / "Disbursement" columns one to
# the right (N->O, O->P, P->Q).
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert()

# The newly inserted column keeps the same width as column M (11 chars)
# but without the "best fit" flag (matches Excel's insert-column behaviour).
$ws.Columns("N").ColumnWidth = 10.14

# Make "Repayment schedule" the active sheet/tab and select Q8 on it
# (this also clears tabSelected on whatever sheet was active before).
$ws.Activate()
[void]$ws.Range("Q8").Select()
